$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 26875
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 100000
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 100000
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -101248
$ws.Range("H64").Value = 3085.7144
$ws.Range("I64").Value = 3023.6365
$ws.Range("J64").Value = 3313.3333
$ws.Range("K64").Value = 3023.6365
$ws.Range("L64").Value = 3313.3333
$ws.Range("M64").Value = -2775.6365
$ws.Range("N64").Value = -3809.3333
$ws.Range("H65").Value = 26875
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 100000
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 500000
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -506240
$ws.Range("H67").Value = 3085.7144
$ws.Range("I67").Value = 3023.6365
$ws.Range("J67").Value = 3313.3333
$ws.Range("K67").Value = 3023.6365
$ws.Range("L67").Value = 3313.3333
$ws.Range("M67").Value = -2165.6365
$ws.Range("N67").Value = -5029.3333
$ws.Range("H138").Value = 2469.1707
$ws.Range("I138").Value = 1439.871
$ws.Range("J138").Value = 3094.8235
$ws.Range("K138").Value = 4319.613
$ws.Range("L138").Value = 9284.470499999999
$ws.Range("M138").Value = 820.3869999999997
$ws.Range("N138").Value = -19564.4705

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 18888
$ws.Range("J9").Value = 18888
$ws.Range("L9").Value = 18888
$ws.Range("N9").Value = -19228
$ws.Range("H20").Value = 18888
$ws.Range("J20").Value = 18888
$ws.Range("L20").Value = 18888
$ws.Range("N20").Value = -19428
$ws.Range("H32").Value = 28333.842
$ws.Range("I32").Value = 5389.7954
$ws.Range("J32").Value = 105990.62
$ws.Range("K32").Value = 5389.7954
$ws.Range("L32").Value = 105990.62
$ws.Range("M32").Value = -5102.7954
$ws.Range("N32").Value = -106564.62
$ws.Range("H74").Value = 686.614
$ws.Range("I74").Value = 702.5909
$ws.Range("J74").Value = 632.53845
$ws.Range("K74").Value = 702.5909
$ws.Range("L74").Value = 632.53845
$ws.Range("M74").Value = 171.4091
$ws.Range("N74").Value = -2380.53845
$ws.Range("H76").Value = 103238.8
$ws.Range("J76").Value = 103238.8
$ws.Range("L76").Value = 103238.8
$ws.Range("N76").Value = -103914.8
$ws.Range("H77").Value = 686.614
$ws.Range("I77").Value = 702.5909
$ws.Range("J77").Value = 632.53845
$ws.Range("K77").Value = 3512.9545
$ws.Range("L77").Value = 3162.69225
$ws.Range("M77").Value = 855.0454999999997
$ws.Range("N77").Value = -11898.69225
$ws.Range("H79").Value = 103238.8
$ws.Range("J79").Value = 103238.8
$ws.Range("L79").Value = 103238.8
$ws.Range("N79").Value = -105578.8
$ws.Range("H92").Value = 32993.332
$ws.Range("J92").Value = 32993.332
$ws.Range("L92").Value = 32993.332
$ws.Range("N92").Value = -37985.332

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 536.6087
$ws.Range("I22").Value = 525.0909
$ws.Range("J22").Value = 790
$ws.Range("K22").Value = 525.0909
$ws.Range("L22").Value = 790
$ws.Range("M22").Value = -352.0909
$ws.Range("N22").Value = -1136
$ws.Range("H82").Value = 17304.223
$ws.Range("I82").Value = 4669.625
$ws.Range("J82").Value = 27411.9
$ws.Range("K82").Value = 4669.625
$ws.Range("L82").Value = 27411.9
$ws.Range("M82").Value = -4286.625
$ws.Range("N82").Value = -28177.9
$ws.Range("H85").Value = 17304.223
$ws.Range("I85").Value = 4669.625
$ws.Range("J85").Value = 27411.9
$ws.Range("K85").Value = 4669.625
$ws.Range("L85").Value = 27411.9
$ws.Range("M85").Value = -3343.625
$ws.Range("N85").Value = -30063.9

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 108.411766
$ws.Range("I7").Value = 63.666668
$ws.Range("J7").Value = 158.75
$ws.Range("K7").Value = 63.666668
$ws.Range("L7").Value = 158.75
$ws.Range("M7").Value = 49.333332
$ws.Range("N7").Value = -384.75
$ws.Range("H99").Value = 1335.5555
$ws.Range("I99").Value = 1252.5
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1252.5
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 245.5
$ws.Range("N99").Value = -4996
$ws.Range("H126").Value = 1335.5555
$ws.Range("I126").Value = 1252.5
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3757.5
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -1287.5
$ws.Range("N126").Value = -10940

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 608.3333
$ws.Range("I47").Value = 412.5
$ws.Range("K47").Value = 1237.5
$ws.Range("M47").Value = -806.5
$ws.Range("H48").Value = 2814.2856
$ws.Range("I48").Value = 300
$ws.Range("J48").Value = 3233.3333
$ws.Range("K48").Value = 900
$ws.Range("L48").Value = 9699.999899999999
$ws.Range("M48").Value = -650
$ws.Range("N48").Value = -10199.9999
$ws.Range("H55").Value = 103490.4
$ws.Range("I55").Value = 3500
$ws.Range("J55").Value = 114600.445
$ws.Range("K55").Value = 10500
$ws.Range("L55").Value = 343801.335
$ws.Range("M55").Value = -10323
$ws.Range("N55").Value = -344155.335
$ws.Range("H64").Value = 11365.272
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 13446.444
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 40339.33199999999
$ws.Range("M64").Value = -5730
$ws.Range("N64").Value = -40879.33199999999
$ws.Range("H67").Value = 11365.272
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 13446.444
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 40339.33199999999
$ws.Range("M67").Value = -5064
$ws.Range("N67").Value = -42211.33199999999
$ws.Range("H68").Value = 1766.8379
$ws.Range("I68").Value = 1193.25
$ws.Range("J68").Value = 1925.069
$ws.Range("K68").Value = 3579.75
$ws.Range("L68").Value = 5775.207
$ws.Range("M68").Value = -2768.75
$ws.Range("N68").Value = -7397.207
$ws.Range("H71").Value = 1766.8379
$ws.Range("I71").Value = 1193.25
$ws.Range("J71").Value = 1925.069
$ws.Range("K71").Value = 10739.25
$ws.Range("L71").Value = 17325.621
$ws.Range("M71").Value = -6683.25
$ws.Range("N71").Value = -25437.621
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("H100").Value = 9600
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 9600
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 28800
$ws.Range("N100").Value = -30422
$ws.Range("H103").Value = 1807
$ws.Range("I103").Value = 262.5
$ws.Range("J103").Value = 3866.3333
$ws.Range("K103").Value = 787.5
$ws.Range("L103").Value = 11598.9999
$ws.Range("M103").Value = 91.5
$ws.Range("N103").Value = -13356.9999
$ws.Range("H107").Value = 681.4783
$ws.Range("I107").Value = 365.8158
$ws.Range("J107").Value = 1068.4193
$ws.Range("K107").Value = 1097.4474
$ws.Range("L107").Value = 3205.2579
$ws.Range("M107").Value = 822.5526
$ws.Range("N107").Value = -7045.257900000001
$ws.Range("H125").Value = 2206.818
$ws.Range("I125").Value = 858.3333
$ws.Range("J125").Value = 2712.5
$ws.Range("K125").Value = 2574.9999
$ws.Range("L125").Value = 8137.5
$ws.Range("M125").Value = 2345.0001
$ws.Range("N125").Value = -17977.5
$ws.Range("H133").Value = 6554.3335
$ws.Range("I133").Value = 2249.75
$ws.Range("J133").Value = 9998
$ws.Range("K133").Value = 6749.25
$ws.Range("L133").Value = 29994
$ws.Range("M133").Value = -1689.25
$ws.Range("N133").Value = -40114
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("M100").ClearContents()

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 29320
$ws.Range("J75").Value = 29320
$ws.Range("L75").Value = 29320
$ws.Range("N75").Value = -31192
$ws.Range("H78").Value = 29320
$ws.Range("J78").Value = 29320
$ws.Range("L78").Value = 87960
$ws.Range("N78").Value = -97320
